$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header cells: "_old" -> "_FV2210", "_new" -> "_FV2304"
$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Temporarily clear the header row's direct formatting so that
# ListObjects.Add doesn't bake the existing look into a new dxf
# (the table's own header styling would otherwise duplicate it).
$headerRange = $ws.Range("A1:U1")
$headerRange.ClearFormats()

# Create a table (ListObject) over the used range, with headers
$range = $ws.Range("A1:U67")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"

# Restore the header row's original formatting (bold, gray fill,
# centered + wrapped text, thin border) directly on the cells.
$headerRange.Font.Bold = $true
$headerRange.Interior.Color = 14277081
$headerRange.HorizontalAlignment = -4108
$headerRange.WrapText = $true
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Freeze the top row (header row)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
